# REVER_DailyTracker_BALAJI.xlsx - "Add files via upload"
#
# Fills in the Application (col C) / Task (col D) columns for the
# 04/03/2021, 05/03/2021, 08/03/2021 and 09/03/2021 rows (sheet rows
# 5, 6, 9 and 10) and moves the sheet's selection to D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (04/03/2021) ---------------------------------------------
$ws.Range("C5").Value = "Hayaai and Mujistore"
# Touching Font.Bold (already False/default) nudges this cell onto the
# plain bordered style used by the other "Application" cells in the
# column (same look, no font/fill change) instead of the unformatted
# blank-cell style it inherited.
$ws.Range("C5").Font.Bold = $false
$ws.Range("D5").Value = "Analyse the mujistore code and support database work for mujistore"

# --- Row 6 (05/03/2021) ---------------------------------------------
$ws.Range("C6").Value = "Hayaai and Mujistore"
$ws.Range("C6").Font.Bold = $false
$ws.Range("D6").Value = "Analyse the mujistore code and support database work for mujistore"

# --- Row 9 (08/03/2021) ---------------------------------------------
$ws.Range("C9").Value = "Mujistore "
$ws.Range("D9").Value = "support for Mujistore issues"

# --- Row 10 (09/03/2021) ---------------------------------------------
$ws.Range("C10").Value = "Mujistore "
$ws.Range("D10").Value = "support for Mujistore issues"

# --- View state: scroll/selection update ------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("D10").Select()
